$d = $word.ActiveDocument

$replacements = @(
    @{old="44÷2="; new="36÷9="},
    @{old="94÷4="; new="37÷9="},
    @{old="57÷3="; new="99÷5="},
    @{old="25÷2="; new="16÷2="},
    @{old="15÷9="; new="57÷7="},
    @{old="58÷2="; new="44÷9="},
    @{old="34÷2="; new="68÷7="},
    @{old="24÷6="; new="36÷9="},
    @{old="68÷8="; new="85÷2="},
    @{old="63÷2="; new="64÷7="},
    @{old="75÷8="; new="60÷4="},
    @{old="87÷3="; new="78÷7="},
    @{old="41÷5="; new="38÷8="},
    @{old="12÷6="; new="10÷7="},
    @{old="62÷3="; new="53÷6="},
    @{old="10÷4="; new="20÷9="},
    @{old="59÷9="; new="32÷3="},
    @{old="18÷3="; new="24÷4="},
    @{old="93÷5="; new="49÷5="},
    @{old="27÷7="; new="66÷4="},
    @{old="72÷4="; new="12÷5="},
    @{old="60÷3="; new="49÷6="},
    @{old="98÷6="; new="77÷8="},
    @{old="66÷5="; new="83÷7="},
    @{old="57÷2="; new="35÷5="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
